$d = $word.ActiveDocument

# Locate the paragraph index (1-based) of "Schattenbahnhof der Realität anpassen"
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Schattenbahnhof der Realität anpassen*") {
        $targetIndex = $i
        break
    }
    $i = $i + 1
}

$target = $d.Paragraphs.Item($targetIndex)

# Collapse to the end of that paragraph (right after the existing bookmark/run,
# before the paragraph mark) and insert the first new list paragraph.
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p1 = $d.Paragraphs.Item($targetIndex + 1)
$p1.Range.Text = "Screen nicht an jeden Baustein übergeben"

# Insert the second new list paragraph after the first one.
$r2 = $p1.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item($targetIndex + 2)
$p2.Range.Text = "Bild vertikales Gleis reparieren"

# Move the (hidden) _GoBack bookmark so it again sits at the end of the last
# edited paragraph, mirroring Word's own behaviour when new text is typed.
$endOfP2 = $p2.Range
$endOfP2.Collapse(0)
$endOfP2.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $endOfP2)
